$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "H 72" record (row 2). This shifts all subsequent rows up by
# one, so the table now ends at row 62 instead of row 63, matching the
# reduced dimension in the target workbook.
$ws.Rows.Item(2).Delete()
